$wb = $excel.ActiveWorkbook

# Sheets 1-4 share the same "Fonte/Tecnologia" row-label structure.
$fonteSheets = @(1, 2, 3, 4)

foreach ($idx in $fonteSheets) {
    $ws = $wb.Worksheets.Item($idx)

    # Grab the header formatting (bold, border, centered) from an existing
    # labeled cell (A2) before we touch it, so the new A1 header cell can
    # reuse the very same style.
    $ws.Range("A2").Copy()
    $ws.Range("A1").PasteSpecial(-4122)
    $ws.Range("A1").Value = "Fonte/Tecnologia"

    # The row labels lose their bold/border/center style and a few of them
    # get corrected spelling/accents.
    $ws.Range("A2").ClearFormats()

    $ws.Range("A3").ClearFormats()
    $ws.Range("A3").Value = "Gás Natural"

    $ws.Range("A4").ClearFormats()
    $ws.Range("A4").Value = "Carvão"

    $ws.Range("A5").ClearFormats()

    $ws.Range("A6").ClearFormats()
    $ws.Range("A6").Value = "Óleos Comb"

    $ws.Range("A7").ClearFormats()

    $ws.Range("A8").ClearFormats()
    $ws.Range("A8").Value = "Eólica"

    $ws.Range("A9").ClearFormats()

    $ws.Range("A10").ClearFormats()

    $ws.Range("A11").ClearFormats()
    $ws.Range("A11").Value = "Pot. Compl."

    $ws.Range("A12").ClearFormats()
}

# Sheet 5: "Emissoes Totais (MtCO2eq)" - header + accents + drop the "Teto" row.
$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("A2").Copy()
$ws5.Range("A1").PasteSpecial(-4122)
$ws5.Range("A1").Value = "Período"

$ws5.Range("A2").ClearFormats()
$ws5.Range("A2").Value = "P.Médio"

$ws5.Range("A3").ClearFormats()
$ws5.Range("A3").Value = "P.Crítico"

$ws5.Rows("4:4").Delete()

# Sheet 6: "Custo Total (bilhões de R$)" - header + accents + updated figures.
$ws6 = $wb.Worksheets.Item(6)

$ws6.Range("A2").Copy()
$ws6.Range("A1").PasteSpecial(-4122)
$ws6.Range("A1").Value = "Tipo Expansão"

# "2015" needs to stay text (like the other sheets' header row), not get
# reinterpreted as a number, while keeping the original bold/border style.
$ws6.Range("B1").NumberFormat = "@"
$ws6.Range("B1").Value = "2015"
$ws6.Range("B1").ClearFormats()
$ws6.Range("A2").Copy()
$ws6.Range("B1").PasteSpecial(-4122)

$ws6.Range("A2").ClearFormats()
$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("B2").Value = 567

$ws6.Range("A3").ClearFormats()
$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("B3").Value = 99
